$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date, insert a new "Jurisdiction" row after "Contact" ---
$meta = $wb.Worksheets.Item("Metadata")

# Update the Date property value (row 8, column B)
$meta.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"

# Insert a new row after "Contact" (row 10) so "Jurisdiction" becomes row 11
$meta.Rows.Item(11).Insert()
$meta.Cells.Item(11, 1).Value = "Jurisdiction"
$meta.Cells.Item(11, 2).Value = ""

# Copy the style used by the other property rows (e.g. row 12, which was row 11 before insert)
$meta.Cells.Item(11, 1).Style = $meta.Cells.Item(12, 1).Style
$meta.Cells.Item(11, 2).Style = $meta.Cells.Item(12, 2).Style
